# Update the "Return_with_prediction" (column G) and "return_pct_change"
# (column H) values for data rows 2-85, plus a one-off correction to the
# "mean_return_pct_change" value in I2 -- per the recorded diff for
# S&P500_returns_compared_quarter.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 2

# New "Return_with_prediction" values (column G), one per row starting at $startRow
$gValues = @(
    0.05827466589841016, 0.06513461836194495, -0.0245898057029326, -0.005005280353070471, -0.02689097679502047, -0.01177810742294736,
    0.005862896202882596, -0.01070679914971268, -0.06048121483629408, -0.05841567825681807, -0.3724713856422338, -0.4165823501648836,
    -0.0243803164075258, -0.02082395709492418, 0.1319814381583955, 0.1515503134163456, 0.1152557440024727, 0.1187911247363713,
    0.08956501910640609, 0.09303514630930471, -0.09913972101773226, -0.1021737610584678, 0.1645349111905927, 0.1723493042693961,
    0.08992923081938585, 0.08589392136060651, -0.1402476163640358, -0.1431779844069648, 0.0491453962518169, 0.04611227313010585,
    0.1163156479122742, 0.1211998253111457, -0.01310417945064889, -0.009902989387291007, 0.03730085483628309, 0.03381165244018244,
    0.09945957467631172, 0.1009795311699805, 0.0312154112119229, 0.03888384873672467, 0.1200430566872892, 0.1249960543188665,
    0.03022879591755538, 0.03967307773274016, 0.06375237553196902, 0.06109720973819019, 0.05117968818501604, 0.04660715296829505,
    0.02287096329384958, 0.02840327117846986, -0.08580245521385596, -0.08041369526160598, 0.04901220076012972, 0.06014393012333171,
    0.04711405073869959, 0.05269221828464761, 0.05711225681876087, 0.05226098250534379, 0.03150696321333957, 0.03491787351517571,
    0.06432817623783894, 0.06520598171726776, 0.02965341724951642, 0.03339832469503045, 0.08076698551228291, 0.08322477414282513,
    -0.02204873618446903, -0.0135694954079274, 0.06754850956833322, 0.08256776905113471, -0.141976198882637, -0.151324181195321,
    0.148909956860133, 0.1568721216574535, -0.003871832189997833, 0.005945211599917738, 0.09680732069173677, 0.09847787780562213,
    -0.223669795498965, -0.2097406981482382, 0.1714130635661462, 0.174552637461663, 0.1091624818887864, 0.1127818040175307
)

# New "return_pct_change" values (column H), one per row starting at $startRow
$hValues = @(
    -9.373760699321917, 15.88338142431029, -21.06282106030316, 56.21436936162208, -140.236192906482, -107.9565043224862,
    203.053756866866, -94.91196277740477, 4.031210858143859, 8.83667637944923, 5.608673204230881, -6.238702419003555,
    -197.5290745901755, 54.04994324237811, -3.486817846540598, 8.661088117801725, -2.050225685422295, -7.676859165718911,
    0.9358086985520511, 6.825182070880617, -6.050877015493134, -0.7183624663165582, 2.137827928964812, 1.027911007054918,
    -0.7992658928256947, -0.07890128100031786, -1.894196263718393, -2.41646967496682, -5.517376117686004, 5.232545309860479,
    7.002245576934087, -2.333182349651732, 16.087033311651, 40.83277322828957, 1.452244209173835, -5.251399662724963,
    -0.8401156224125716, 3.675067403668406, -7.338626337255988, 20.68711190546636, -0.7132630806288016, -2.182872112776352,
    -23.79289873092075, 27.31091846728669, 12.60156263919403, 4.136722948414437, 3.902699855298778, 2.281011789132631,
    -13.65043579713264, 1.382797463602878, 1.28985150983415, -0.2462013007160604, -2.02791934728214, 6.855946439009461,
    -4.688964987638402, 38.7222225155595, -0.8709607608757222, -8.335469097575427, 14.73766848455711, 30.78920511616709,
    3.012055022520502, 2.071473181155906, 6.899003100964459, -5.727023994294394, 3.969719442670271, 5.520343601962471,
    -1.393579200924948, 29.108854875223, -6.164200978850482, 3.962087741420979, 7.592535096485426, 1.144123659759551,
    -1.003354007620764, 4.270462860117357, -273.5625099069332, 369.2714566508805, 7.613609917333086, 1.626925570036043,
    -3.343489825467979, 1.582809675302308, 2.264317003184549, -0.8386659537782877, 2.85762134341584, 7.864471062878129
)

for ($i = 0; $i -lt $gValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $gValues[$i]
    $ws.Cells.Item($row, 8).Value = $hValues[$i]
}

# One-off correction to I2 ("mean_return_pct_change")
$ws.Cells.Item(2, 9).Value = 1.125878855932595
